# Update F3 product_type from "rulebook" to "box set"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "box set"

# Move the active selection to B12 (matches the saved state in the diff)
$ws.Range("B12").Select()
